$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.353.82'
$ws.Range('E2').Value = '  +1.90%  '

$ws.Range('D3').Value = '2.396.20'
$ws.Range('E3').Value = '  +7.68%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +11.11%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.77%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.654'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.45%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.654'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.14%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.26%  '

$ws.Range('E11').Value = '  +3.52%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.62'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.33%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.61'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +17.99%  '

$ws.Range('E14').Value = '  -0.14%  '

$ws.Range('E15').Value = '  +2.29%  '

$ws.Range('D16').Value = '2.752.58'
$ws.Range('E16').Value = '  +7.73%  '

$ws.Range('D17').Value = '2.389.98'
$ws.Range('E17').Value = '  +7.27%  '

$ws.Range('D18').Value = '43.346.58'
$ws.Range('E18').Value = '  +2.31%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.44'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.96%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000109'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.11%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '76.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.04%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '271.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +15.26%  '

$ws.Range('E23').Value = '  +3.35%  '

$ws.Range('E24').Value = '  +1.89%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.80%  '

$ws.Range('E26').Value = '  +3.81%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.05%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.06%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '177.11'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.14%  '

$ws.Range('E30').Value = '  -0.30%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.36%  '

$ws.Range('E32').Value = '  +1.91%  '

$ws.Range('E33').Value = '  +5.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.90'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.07%  '

$ws.Range('E35').Value = '  +6.31%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.87'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.66%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.08'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.94%  '

$ws.Range('E38').Value = '  -1.90%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.109'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.54%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.84'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +18.74%  '

$ws.Range('E41').Value = '  +22.50%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '127.41'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +25.81%  '

$ws.Range('E43').Value = '  +1.53%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.00%  '

$ws.Range('E45').Value = '  +0.33%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.62'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.09%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.62'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.48%  '

$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +61.68%  '

$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.79%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.32'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.82%  '

$ws.Range('D51').Value = '1.606.90'
$ws.Range('E51').Value = '  +12.70%  '
